$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80: Cleansing the Wicked Humours
$ws.Range("H80").Value = 30305008
$ws.Range("I80").Value = 47620064
$ws.Range("J80").Value = 3662.5
$ws.Range("K80").Value = 142860192
$ws.Range("L80").Value = 10987.5
$ws.Range("M80").Value = -142859194
$ws.Range("N80").Value = -12983.5

# Row 83: Washing Away the Sins (L)
$ws.Range("H83").Value = 30305008
$ws.Range("I83").Value = 47620064
$ws.Range("J83").Value = 3662.5
$ws.Range("K83").Value = 428580576
$ws.Range("L83").Value = 32962.5
$ws.Range("M83").Value = -428575584
$ws.Range("N83").Value = -42946.5

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1575825.5
$ws.Range("I137").Value = 3853697.8
$ws.Range("J137").Value = 4879.0347
$ws.Range("K137").Value = 11561093.4
$ws.Range("L137").Value = 14637.1041
$ws.Range("M137").Value = -11558543.4
$ws.Range("N137").Value = -19737.1041

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2598.9836
$ws.Range("I138").Value = 1624.5714
$ws.Range("J138").Value = 3425.7576
$ws.Range("K138").Value = 4873.7142
$ws.Range("L138").Value = 10277.2728
$ws.Range("M138").Value = 266.2857999999997
$ws.Range("N138").Value = -20557.2728

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 4044.4443
$ws.Range("I141").Value = 2023.1666
$ws.Range("J141").Value = 8087
$ws.Range("K141").Value = 6069.4998
$ws.Range("L141").Value = 24261
$ws.Range("M141").Value = -889.4997999999996
$ws.Range("N141").Value = -34621

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 1464.7297
$ws.Range("I2").Value = 1637.7693
$ws.Range("J2").Value = 1055.7273
$ws.Range("K2").Value = 1637.7693
$ws.Range("L2").Value = 1055.7273
$ws.Range("M2").Value = -1524.7693
$ws.Range("N2").Value = -1281.7273

# Row 7: Distill It Yourself
$ws.Range("H7").Value = 28888.889
$ws.Range("J7").Value = 28888.889
$ws.Range("L7").Value = 28888.889
$ws.Range("N7").Value = -29116.889

# Row 32: Ingot We Trust
$ws.Range("H32").Value = 9739.769
$ws.Range("I32").Value = 9688.925999999999
$ws.Range("J32").Value = 9922.799999999999
$ws.Range("K32").Value = 9688.925999999999
$ws.Range("L32").Value = 9922.799999999999
$ws.Range("M32").Value = -9401.925999999999
$ws.Range("N32").Value = -10496.8

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 2018.4878
$ws.Range("I74").Value = 1627
$ws.Range("J74").Value = 4302.1665
$ws.Range("K74").Value = 1627
$ws.Range("L74").Value = 4302.1665
$ws.Range("M74").Value = -753
$ws.Range("N74").Value = -6050.1665

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 2018.4878
$ws.Range("I77").Value = 1627
$ws.Range("J77").Value = 4302.1665
$ws.Range("K77").Value = 8135
$ws.Range("L77").Value = 21510.8325
$ws.Range("M77").Value = -3767
$ws.Range("N77").Value = -30246.8325

# Row 88: The Mast Chance
$ws.Range("H88").Value = 17871350
$ws.Range("I88").Value = 66672000
$ws.Range("J88").Value = 3231155.5
$ws.Range("K88").Value = 66672000
$ws.Range("L88").Value = 3231155.5
$ws.Range("M88").Value = -66671594
$ws.Range("N88").Value = -3231967.5

# Row 91: The Rose and the Riveter (L)
$ws.Range("H91").Value = 17871350
$ws.Range("I91").Value = 66672000
$ws.Range("J91").Value = 3231155.5
$ws.Range("K91").Value = 66672000
$ws.Range("L91").Value = 3231155.5
$ws.Range("M91").Value = -66670596
$ws.Range("N91").Value = -3233963.5

# Row 116: No Scope
$ws.Range("H116").Value = 1464.7297
$ws.Range("I116").Value = 1637.7693
$ws.Range("J116").Value = 1055.7273
$ws.Range("K116").Value = 1637.7693
$ws.Range("L116").Value = 1055.7273
$ws.Range("M116").Value = 656.2307000000001
$ws.Range("N116").Value = -5643.7273

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 1464.7297
$ws.Range("I3").Value = 1637.7693
$ws.Range("J3").Value = 1055.7273
$ws.Range("K3").Value = 1637.7693
$ws.Range("L3").Value = 1055.7273
$ws.Range("M3").Value = -1523.7693
$ws.Range("N3").Value = -1283.7273

# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 2763.1667
$ws.Range("I86").Value = 2115.9
$ws.Range("J86").Value = 5999.5
$ws.Range("K86").Value = 2115.9
$ws.Range("L86").Value = 5999.5
$ws.Range("M86").Value = -992.9000000000001
$ws.Range("N86").Value = -8245.5

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 2763.1667
$ws.Range("I89").Value = 2115.9
$ws.Range("J89").Value = 5999.5
$ws.Range("K89").Value = 10579.5
$ws.Range("L89").Value = 29997.5
$ws.Range("M89").Value = -4963.5
$ws.Range("N89").Value = -41229.5

# Row 108: Fire Sale
$ws.Range("H108").Value = 39000
$ws.Range("J108").Value = 39000
$ws.Range("L108").Value = 39000
$ws.Range("N108").Value = -46680

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2537
$ws.Range("I134").Value = 2223.72
$ws.Range("J134").Value = 2997.7058
$ws.Range("K134").Value = 6671.16
$ws.Range("L134").Value = 8993.117400000001
$ws.Range("M134").Value = -4136.16
$ws.Range("N134").Value = -14063.1174

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 3626827.2
$ws.Range("I31").Value = 2157.8386
$ws.Range("J31").Value = 5468872.5
$ws.Range("K31").Value = 2157.8386
$ws.Range("L31").Value = 5468872.5
$ws.Range("M31").Value = -1862.8386
$ws.Range("N31").Value = -5469462.5

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 3626827.2
$ws.Range("I34").Value = 2157.8386
$ws.Range("J34").Value = 5468872.5
$ws.Range("K34").Value = 2157.8386
$ws.Range("L34").Value = 5468872.5
$ws.Range("M34").Value = -1955.8386
$ws.Range("N34").Value = -5469276.5

# Row 55: Ready for a Rematch
$ws.Range("H55").Value = 3340
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 4233.3335
$ws.Range("K55").Value = 2000
$ws.Range("L55").Value = 4233.3335
$ws.Range("M55").Value = -1685
$ws.Range("N55").Value = -4863.3335

$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face
$ws.Range("H68").Value = 1487.3784
$ws.Range("J68").Value = 1634.3208
$ws.Range("L68").Value = 4902.9624
$ws.Range("N68").Value = -6524.9624

# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 1487.3784
$ws.Range("J71").Value = 1634.3208
$ws.Range("L71").Value = 14708.8872
$ws.Range("N71").Value = -22820.8872

# Row 133: Friends Are Food
$ws.Range("H133").Value = 6745.385
$ws.Range("I133").Value = 8700
$ws.Range("J133").Value = 6582.5
$ws.Range("K133").Value = 26100
$ws.Range("L133").Value = 19747.5
$ws.Range("M133").Value = -21040
$ws.Range("N133").Value = -29867.5

# Row 134: Don't Knock It Till You've Tried It
$ws.Range("H134").Value = 43659124
$ws.Range("I134").Value = 47816564
$ws.Range("J134").Value = 6025
$ws.Range("K134").Value = 143449692
$ws.Range("L134").Value = 18075
$ws.Range("M134").Value = -143444622
$ws.Range("N134").Value = -28215

$ws = $wb.Worksheets.Item("GSM")
# Row 110: Slimming Down
$ws.Range("H110").Value = 42000
$ws.Range("J110").Value = 42000
$ws.Range("L110").Value = 42000
$ws.Range("N110").Value = -50180
